# Weekly driver report update for 2025-05-05
# Applies to the "Driver Summary" sheet of IPO_driver_summary.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the obsolete "Bad Driver" row (Intel(R) Wi-Fi 6E AX211 160MHz -
#    23.30.0.6), which also shifts every row below it up by one. This moves
#    the old Totals row (5 -> 4), the blank spacer rows (6-10 -> 5-9), the
#    "Good Drivers" header (11 -> 10), the good-driver column headers
#    (12 -> 11) and the thirteen good-driver data rows (13-25 -> 12-24) while
#    preserving their existing cell styles.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# 2. Refresh the remaining "Bad Driver" row values (row 3).
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 3).Value = 657     # Critical Minutes
$ws.Cells.Item(3, 4).Value = 91.5    # Good Roaming Calculation (%)

# ---------------------------------------------------------------------------
# 3. Refresh the Totals row (now row 4 after the shift).
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = 3       # Client Count total
$ws.Cells.Item(4, 3).Value = 657     # Critical Minutes total

# ---------------------------------------------------------------------------
# 4. Refresh the "Good Drivers" data table (now rows 12-24 after the shift)
#    with this week's numbers. Rows are written in the new report order.
# ---------------------------------------------------------------------------
$adapterName = "intel(r) wi-fi 6e ax211 160mhz"

$goodDrivers = @(
    @("22.180.0.4", 17609,   69,   0,    57,   17678,   99.59999999999999, "2022-10-17"),
    @("22.250.0.4", 1293197, 4322, 1009, 1990, 1298528, 99.59999999999999, "2023-07-25"),
    @("22.220.0.4", 31517,   112,  0,    66,   31629,   99.59999999999999, "2023-03-28"),
    @("23.10.0.8",  467311,  772,  931,  706,  469014,  99.59999999999999, "2023-10-30"),
    @("23.120.0.3", 455081,  1861, 52,   639,  456994,  99.59999999999999, "2025-02-05"),
    @("22.230.0.8", 1787924, 3326, 2614, 3038, 1793864, 99.7,               "2023-05-08"),
    @("23.70.2.3",  218767,  334,  313,  573,  219414,  99.7,               "2024-07-23"),
    @("22.110.1.1", 135467,  189,  263,  196,  135919,  99.7,               "2022-01-01"),
    @("23.100.0.4", 240434,  421,  37,   409,  240892,  99.8,               "2024-11-10"),
    @("23.80.1.3",  151287,  285,  75,   332,  151647,  99.8,               "2024-09-03"),
    @("22.100.1.1", 272039,  213,  131,  316,  272383,  99.90000000000001, "2022-05-01"),
    @("22.150.0.3", 14561,   0,    0,    59,   14561,   100,                "2022-05-23"),
    @("22.150.3.1", 12018,   0,    0,    61,   12018,   100,                "2022-08-29")
)

$startRow = 12
$endRow = $startRow + $goodDrivers.Count - 1

# Pre-format the "driver vintage" column as text so that date-like strings
# such as "2022-10-17" are stored as literal text (matching the source
# report) instead of being auto-converted into Excel date serial numbers.
$vintageRange = $ws.Range("J" + $startRow + ":J" + $endRow)
$vintageRange.NumberFormat = "@"

for ($i = 0; $i -lt $goodDrivers.Count; $i++) {
    $row = $startRow + $i
    $driver = $goodDrivers[$i]

    $ws.Cells.Item($row, 1).Value = "Intel(R) Wi-Fi 6E AX211 160MHz - " + $driver[0]
    $ws.Cells.Item($row, 2).Value = $driver[1]
    $ws.Cells.Item($row, 3).Value = $driver[2]
    $ws.Cells.Item($row, 4).Value = $driver[3]
    $ws.Cells.Item($row, 5).Value = $driver[4]
    $ws.Cells.Item($row, 6).Value = $driver[5]
    $ws.Cells.Item($row, 7).Value = $adapterName
    $ws.Cells.Item($row, 8).Value = $driver[0]
    $ws.Cells.Item($row, 9).Value = $driver[6]
    $ws.Cells.Item($row, 10).Value = $driver[7]
}

# Restore the default (unstyled) look of the vintage column now that the
# values have safely been stored as text.
$vintageRange.Style = "Normal"
